$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4678.2
$ws.Range("I74").Value = 4554.2856
$ws.Range("J74").Value = 4967.3335
$ws.Range("K74").Value = 4554.2856
$ws.Range("L74").Value = 4967.3335
$ws.Range("M74").Value = -3618.2856
$ws.Range("N74").Value = -6839.3335

$ws.Range("H77").Value = 4678.2
$ws.Range("I77").Value = 4554.2856
$ws.Range("J77").Value = 4967.3335
$ws.Range("K77").Value = 22771.428
$ws.Range("L77").Value = 24836.6675
$ws.Range("M77").Value = -18091.428
$ws.Range("N77").Value = -34196.6675

$ws.Range("H100").Value = 2006.0476
$ws.Range("J100").Value = 2527.25
$ws.Range("L100").Value = 2527.25
$ws.Range("N100").Value = -3609.25

$ws.Range("H137").Value = 2352.8696
$ws.Range("I137").Value = 1752.0513
$ws.Range("J137").Value = 5700.2856
$ws.Range("K137").Value = 5256.1539
$ws.Range("L137").Value = 17100.8568
$ws.Range("M137").Value = -2706.1539
$ws.Range("N137").Value = -22200.8568

$ws.Range("H140").Value = 86028.8
$ws.Range("J140").Value = 86028.8
$ws.Range("L140").Value = 86028.8
$ws.Range("N140").Value = -96388.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8217.383
$ws.Range("I61").Value = 5275.758
$ws.Range("J61").Value = 15151.214
$ws.Range("K61").Value = 5275.758
$ws.Range("L61").Value = 15151.214
$ws.Range("M61").Value = -5063.758
$ws.Range("N61").Value = -15575.214

$ws.Range("H74").Value = 81209.92999999999
$ws.Range("I74").Value = 93745.87
$ws.Range("J74").Value = 1815.6666
$ws.Range("K74").Value = 93745.87
$ws.Range("L74").Value = 1815.6666
$ws.Range("M74").Value = -92871.87
$ws.Range("N74").Value = -3563.6666

$ws.Range("H77").Value = 81209.92999999999
$ws.Range("I77").Value = 93745.87
$ws.Range("J77").Value = 1815.6666
$ws.Range("K77").Value = 468729.35
$ws.Range("L77").Value = 9078.333000000001
$ws.Range("M77").Value = -464361.35
$ws.Range("N77").Value = -17814.333

$ws.Range("H97").Value = 789.7778
$ws.Range("I97").Value = 718.58826
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 718.58826
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -222.58826
$ws.Range("N97").Value = -2992

$ws.Range("H102").Value = 3716.1538
$ws.Range("I102").Value = 3511
$ws.Range("K102").Value = 3511
$ws.Range("M102").Value = -1889

$ws.Range("H132").Value = 7490.6
$ws.Range("I132").Value = 9683.532999999999
$ws.Range("J132").Value = 4201.2
$ws.Range("K132").Value = 29050.599
$ws.Range("L132").Value = 12603.6
$ws.Range("M132").Value = -26520.599
$ws.Range("N132").Value = -17663.6

$ws.Range("H136").Value = 8217.383
$ws.Range("I136").Value = 5275.758
$ws.Range("J136").Value = 15151.214
$ws.Range("K136").Value = 15827.274
$ws.Range("L136").Value = 45453.642
$ws.Range("M136").Value = -13277.274
$ws.Range("N136").Value = -50553.642

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 834.375
$ws.Range("I94").Value = 665.75
$ws.Range("J94").Value = 1677.5
$ws.Range("K94").Value = 665.75
$ws.Range("L94").Value = 1677.5
$ws.Range("M94").Value = -214.75
$ws.Range("N94").Value = -2579.5

$ws.Range("H134").Value = 34757.71
$ws.Range("I134").Value = 2565.55
$ws.Range("J134").Value = 93288.91
$ws.Range("K134").Value = 7696.650000000001
$ws.Range("L134").Value = 279866.73
$ws.Range("M134").Value = -5161.650000000001
$ws.Range("N134").Value = -284936.73

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2832.7715
$ws.Range("I31").Value = 2136.3125
$ws.Range("J31").Value = 3419.2632
$ws.Range("K31").Value = 2136.3125
$ws.Range("L31").Value = 3419.2632
$ws.Range("M31").Value = -1841.3125
$ws.Range("N31").Value = -4009.2632

$ws.Range("H34").Value = 2832.7715
$ws.Range("I34").Value = 2136.3125
$ws.Range("J34").Value = 3419.2632
$ws.Range("K34").Value = 2136.3125
$ws.Range("L34").Value = 3419.2632
$ws.Range("M34").Value = -1934.3125
$ws.Range("N34").Value = -3823.2632

$ws.Range("H58").Value = 2757207.2
$ws.Range("I58").Value = 5682927.5
$ws.Range("J58").Value = 3588.4707
$ws.Range("K58").Value = 5682927.5
$ws.Range("L58").Value = 3588.4707
$ws.Range("M58").Value = -5682724.5
$ws.Range("N58").Value = -3994.4707

$ws.Range("H132").Value = 2424.83
$ws.Range("I132").Value = 2140.2104
$ws.Range("J132").Value = 3145.8667
$ws.Range("K132").Value = 6420.6312
$ws.Range("L132").Value = 9437.6001
$ws.Range("M132").Value = -3890.6312
$ws.Range("N132").Value = -14497.6001

$ws.Range("H134").Value = 2561.9678
$ws.Range("I134").Value = 2325.7917
$ws.Range("J134").Value = 3371.7144
$ws.Range("K134").Value = 6977.375100000001
$ws.Range("L134").Value = 10115.1432
$ws.Range("M134").Value = -4442.375100000001
$ws.Range("N134").Value = -15185.1432

$ws.Range("H136").Value = 2757207.2
$ws.Range("I136").Value = 5682927.5
$ws.Range("J136").Value = 3588.4707
$ws.Range("K136").Value = 17048782.5
$ws.Range("L136").Value = 10765.4121
$ws.Range("M136").Value = -17046232.5
$ws.Range("N136").Value = -15865.4121

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 2321.1667
$ws.Range("J94").Value = 2756.75
$ws.Range("L94").Value = 8270.25
$ws.Range("N94").Value = -9622.25

$ws.Range("H131").Value = 18413.127
$ws.Range("I131").Value = 1420.6666
$ws.Range("J131").Value = 24785.3
$ws.Range("K131").Value = 4261.9998
$ws.Range("L131").Value = 74355.89999999999
$ws.Range("M131").Value = 778.0002000000004
$ws.Range("N131").Value = -84435.89999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 7568.2383
$ws.Range("I132").Value = 2829.0715
$ws.Range("J132").Value = 17046.572
$ws.Range("K132").Value = 8487.2145
$ws.Range("L132").Value = 51139.716
$ws.Range("M132").Value = -5957.2145
$ws.Range("N132").Value = -56199.716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2066.25
$ws.Range("I93").Value = 1640
$ws.Range("K93").Value = 1640
$ws.Range("M93").Value = -392

$ws.Range("H100").Value = 4361.769
$ws.Range("I100").Value = 3882.0908
$ws.Range("J100").Value = 7000
$ws.Range("K100").Value = 3882.0908
$ws.Range("L100").Value = 7000
$ws.Range("M100").Value = -3341.0908
$ws.Range("N100").Value = -8082

$ws.Range("H132").Value = 3229.6155
$ws.Range("I132").Value = 3144.125
$ws.Range("J132").Value = 3366.4
$ws.Range("K132").Value = 9432.375
$ws.Range("L132").Value = 10099.2
$ws.Range("M132").Value = -6902.375
$ws.Range("N132").Value = -15159.2

$ws.Range("H136").Value = 4209.2
$ws.Range("I136").Value = 2623.1562
$ws.Range("J136").Value = 6415.8696
$ws.Range("K136").Value = 7869.4686
$ws.Range("L136").Value = 19247.6088
$ws.Range("M136").Value = -5319.4686
$ws.Range("N136").Value = -24347.6088

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 25000
$ws.Range("J21").Value = 25000
$ws.Range("L21").Value = 25000
$ws.Range("N21").Value = -25470

$ws.Range("H35").Value = 25000
$ws.Range("J35").Value = 25000
$ws.Range("L35").Value = 25000
$ws.Range("N35").Value = -25580

$ws.Range("H64").Value = 38114
$ws.Range("J64").Value = 38114
$ws.Range("L64").Value = 38114
$ws.Range("N64").Value = -38610

$ws.Range("H67").Value = 38114
$ws.Range("J67").Value = 38114
$ws.Range("L67").Value = 38114
$ws.Range("N67").Value = -39830

$ws.Range("H100").Value = 2130.8333
$ws.Range("I100").Value = 427.33334
$ws.Range("J100").Value = 3834.3333
$ws.Range("K100").Value = 854.66668
$ws.Range("L100").Value = 7668.6666
$ws.Range("M100").Value = -313.66668
$ws.Range("N100").Value = -8750.6666

$ws.Range("H122").Value = 14426.75
$ws.Range("I122").Value = 851
$ws.Range("K122").Value = 2553
$ws.Range("M122").Value = -103

$ws.Range("H132").Value = 2269.258
$ws.Range("I132").Value = 1460.875
$ws.Range("K132").Value = 4382.625
$ws.Range("M132").Value = -1852.625

$ws.Range("H136").Value = 5320.9375
$ws.Range("I136").Value = 7236.727
$ws.Range("J136").Value = 4317.4287
$ws.Range("K136").Value = 21710.181
$ws.Range("L136").Value = 12952.2861
$ws.Range("M136").Value = -19160.181
$ws.Range("N136").Value = -18052.2861
